# Insert a new data row for "Femacal de La Calera" (Haba) above the
# current row 32, shifting all subsequent rows down by one, and fill
# the new row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 32 (pushes old rows 32..72 down to 33..73)
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record's data
$ws.Cells.Item(32, 1).Value = 3
$ws.Cells.Item(32, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44483
$ws.Cells.Item(32, 5).Value = 5
$ws.Cells.Item(32, 6).Value = 100112026
$ws.Cells.Item(32, 7).Value = "Haba"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 90
$ws.Cells.Item(32, 11).Value = 8000
$ws.Cells.Item(32, 12).Value = 8500
$ws.Cells.Item(32, 13).Value = 8278
$ws.Cells.Item(32, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(32, 16).Value = 331
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
